$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.730.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.101.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5196"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4397"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.78"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09424"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.80"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.122.78"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.820"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.209"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.242"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.779.73"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.61"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.368.98"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.41"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.512"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.63"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.136"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.715"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1053"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.212"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.954"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.346"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.49"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02581"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06738"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7011"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.336"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2220"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6833"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.38"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.352"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.621"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000356"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.17%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.53"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.57%  "
